# "fixed: diff = buy - sell"
# Refresh the previously-wrong day-over-day % change for 2020-08-21 (C10)
# and append five new days of data (rows 11-15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix existing row 10 (C10) ---
$ws.Range("C10").Value = 0.07205276134393998

# --- append new rows 11-15 ---
$data = @(
    @{ Row = 11; Date = "2020-08-24"; B = 120516.3; C = -0.3117295125876497; D = 0.8564025016361827; E = 20207.6;  F = -0.2834845014289463; G = 0.1435974983638174 },
    @{ Row = 12; Date = "2020-08-25"; B = 186783.1; C = 0.5498575711335313;  D = 0.8202148560146774; E = 40941.5;  F = 1.026044656465884;   G = 0.1797851439853226 },
    @{ Row = 13; Date = "2020-08-26"; B = 126490.4; C = -0.3227952635971884; D = 0.8414691503310584; E = 23830.5;  F = -0.4179377892847111; G = 0.1585308496689415 },
    @{ Row = 14; Date = "2020-08-27"; B = 111436.7; C = -0.1190106126630955; D = 0.8817938674579624; E = 14938.3;  F = -0.3731436604351567; G = 0.1182061325420376 },
    @{ Row = 15; Date = "2020-08-28"; B = 353768;   C = 2.174609441952247;   D = 0.8631835272703317; E = 56073;   F = 2.753639972419887;   G = 0.1368164727296683 }
)

$formatSource = $ws.Range("A10")

foreach ($item in $data) {
    $r = $item.Row

    $cellA = $ws.Cells.Item($r, 1)
    # force text (avoid Excel auto-converting "2020-08-24" into a date serial)
    $cellA.NumberFormat = "@"
    $cellA.Value = $item.Date

    # copy the exact cell style used by the other date cells (A4:A10)
    $formatSource.Copy()
    $cellA.PasteSpecial(-4122)  # xlPasteFormats

    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}
